$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Shared-string order matters: first-use order below reproduces the
# --- sharedStrings.xml index assignment (9..16) from the target diff.

# idx 9 "Win32 (low frag)" -- first used at A3
$ws.Range("A3").Value = "Win32 (low frag)"

# idx 10 "noinline, noforceinline" -- first used at A37
$ws.Range("A37").Value = "noinline, noforceinline"
$ws.Range("B37").Value = 748133

# idx 11 "noforceinline" -- first used at A38
$ws.Range("A38").Value = "noforceinline"
$ws.Range("B38").Value = 750395

# idx 12 "all enabled" -- first used at A39
$ws.Range("A39").Value = "all enabled"
$ws.Range("B39").Value = 749212

# idx 13 "noforceinline + LTCG" -- first used at A40
$ws.Range("A40").Value = "noforceinline + LTCG"
$ws.Range("B40").Value = 742156

# idx 14 "One thread, 10000 record SpeedTest x64:" -- first used at A36
$ws.Range("A36").Value = "One thread, 10000 record SpeedTest x64:"

# idx 15 "Removed cold code from GetThreadCache" -- first used at A41
$ws.Range("A41").Value = "Removed cold code from GetThreadCache"

# idx 16 "noforceinline, /O1" -- first used at A44
$ws.Range("A44").Value = "noforceinline, /O1"
$ws.Range("B44").Value = 746683

# --- Remaining rows that reuse existing shared strings ---

# Row 21: filler series 1..5 (same pattern as rows 2/11/31)
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 5

# Row 22: reuse "Win32 (low frag)" (idx 9)
$ws.Range("A22").Value = "Win32 (low frag)"

# Row 23
$ws.Range("A23").Value = "nedmalloc v1.06"
$ws.Range("F23").Value = 597722
$ws.Range("G23").Formula = "=F23/F$22"

# Row 24
$ws.Range("A24").Value = "nedmalloc v1.06 (sysalloc)"
$ws.Range("G24").Formula = "=F24/F$22"

# Row 42/43: reuse "noforceinline" / "all enabled"
$ws.Range("A42").Value = "noforceinline"
$ws.Range("B42").Value = 753032

$ws.Range("A43").Value = "all enabled"
$ws.Range("B43").Value = 747830

# --- Edits to existing rows 33 / 34 ---
$ws.Range("B33").Value = 626096
$ws.Range("F33").Value = 542544
$ws.Range("G33").Formula = "=F33/F$32"

$ws.Range("F34").ClearContents()
$ws.Range("G34").Formula = "=F34/F$32"

# --- View / selection state ---
$ws.Range("J24").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

# --- Page setup (portrait orientation) ---
$ws.PageSetup.Orientation = 1
$null = 0
